# Atualização automática do arquivo Excel
# Adds three new records (rows 111-113) to the "Balanço Geral" data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Record", "Balanço Geral", "Infraestrutura", "2025-04-10T13:21", "Negativo", "Moradores e turistas reclamam da falta de infraestrutura do Morro do Itaoca. Banheiro sem porta e sujo. Entrevista com ciclista, falando que o local precisa de mais atenção, inclusive na Ponte do Marimbondo; com o instrutor de parapente, ressaltando que o morro atende a várias tribos; com o piloto de parapente. *sem nota da prefeitura* *matéria* Apresentadora questionou sobre o ICMS Verde `n"),
    @("Record", "Balanço Geral", "Saúde", "2025-04-10T13:27", "Positivo", "Blitz da Saúde promovido pelo Hospital Geral Prontocardio - parceria com a Prefeitura de Campos. Entrevista com a coord. de atendimento Prontocardio, Jeniffer Ferreira. "),
    @("Record", "Balanço Geral", "Social", "2025-04-10T13:45", "Neutro", "ONG Nova Vida oferece aulas de balé para crianças e adolescentes da comunidade. Atualmente, sobrevive da ajuda da comunidade. Em breve, vai atender crianças atípicas. Entrevista com profissionais, alunos, mães e profissional que já foi assistida. Entrevista com psicóloga sobre a importância do esporte na vida das crianças. Na entrevista do fundador da ONG, Sérgio Bessa, ele falou da expectativa de que a prefeitura torne a ONG Utilidade Pública para firmar parceria. ")
)

$startRow = 111
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}
